$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data between row 2 and row 3 for columns D, J, K, L, M, P
# Row 2 current -> D2=44203 J2=27 K2=7000 L2=8000 M2=7556 P2=756
# Row 3 current -> D3=44211 J3=28 K3=8000 L3=8500 M3=8214 P3=821
# Target: row2 gets row3's old values, row3 gets row2's old values

$ws.Range("D2").Value2 = 44211
$ws.Range("J2").Value2 = 28
$ws.Range("K2").Value2 = 8000
$ws.Range("L2").Value2 = 8500
$ws.Range("M2").Value2 = 8214
$ws.Range("P2").Value2 = 821

$ws.Range("D3").Value2 = 44203
$ws.Range("J3").Value2 = 27
$ws.Range("K3").Value2 = 7000
$ws.Range("L3").Value2 = 8000
$ws.Range("M3").Value2 = 7556
$ws.Range("P3").Value2 = 756
